# Storage of Struct.xlsx - update note of C
# Appends a sentence to the big explanatory note (L5, merged L5:P21),
# adds a new "index" offset column (G34:G45 with label in F45), and
# moves the active selection/viewport as the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Extend the explanatory paragraph in L5 with the extra sentence about
#    using 4*index as a subscript for fast access.
$ws.Range("L5").Value = $ws.Range("L5").Value() + "使用4*index (0,1,...)这样的下标时访问速度会很高。"

# 2. Add the new "index" column next to the second (no-padding) memory
#    table: row 45 (offset 0x10, variable a) is index 0, and each row
#    above it decreases by 1 up to row 34 (offset 0x5) = -11. Label the
#    column with the new "index" header text in F45.
$ws.Range("F45").Value = "index"
$ws.Range("G45").Value = 0
$ws.Range("G44").Value = -1
$ws.Range("G43").Value = -2
$ws.Range("G42").Value = -3
$ws.Range("G41").Value = -4
$ws.Range("G40").Value = -5
$ws.Range("G39").Value = -6
$ws.Range("G38").Value = -7
$ws.Range("G37").Value = -8
$ws.Range("G36").Value = -9
$ws.Range("G35").Value = -10
$ws.Range("G34").Value = -11

# 3. Leave the view where the author left it: scrolled up a bit with the
#    new G45 cell selected instead of the old K31:N41 merge.
$ws.Range("G45").Select()
$excel.ActiveWindow.ScrollRow = 10
